$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.1
$ws.Range("I2").Value = 3.8
$ws.Range("AH2").Value = 17

# Row 3 updates
$ws.Range("H3").Value = 2.9
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.25
$ws.Range("Q3").Value = 2.88
$ws.Range("R3").Value = 1.4
$ws.Range("AC3").Value = 5.5
$ws.Range("AO3").Value = 13
$ws.Range("AX3").Value = 26
